$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per the refreshed cryptos feed.
# D-column values that Excel would auto-parse as a plain number (single dot)
# are written with a leading apostrophe so they stay literal text, exactly
# matching the original formatting (e.g. trailing zeros, leading zeros kept).

$ws.Range("D2").Value = "28.155.41"
$ws.Range("E2").Value = "  +0.24%  "

$ws.Range("D3").Value = "1.884.71"
$ws.Range("E3").Value = "  -0.75%  "

$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.38%  "

$ws.Range("D5").Value = "'313.76"
$ws.Range("E5").Value = "  -0.15%  "

$ws.Range("E6").Value = "  +0.34%  "

$ws.Range("D7").Value = "'0.5040"
$ws.Range("E7").Value = "  +0.26%  "

$ws.Range("D8").Value = "'0.3828"
$ws.Range("E8").Value = "  -1.86%  "

$ws.Range("D9").Value = "'0.08556"
$ws.Range("E9").Value = "  -7.42%  "

$ws.Range("D10").Value = "'1.117"
$ws.Range("E10").Value = "  -0.98%  "

$ws.Range("D11").Value = "'41.36"
$ws.Range("E11").Value = "  -1.08%  "

$ws.Range("D12").Value = "'6.268"
$ws.Range("E12").Value = "  -2.10%  "

$ws.Range("D13").Value = "1.881.51"
$ws.Range("E13").Value = "  -0.46%  "

$ws.Range("D14").Value = "'20.66"
$ws.Range("E14").Value = "  -0.75%  "

$ws.Range("D15").Value = "'7.218"
$ws.Range("E15").Value = "  -1.23%  "

$ws.Range("E16").Value = "  +0.38%  "

$ws.Range("D17").Value = "'0.00001099"
$ws.Range("E17").Value = "  -1.02%  "

$ws.Range("D18").Value = "'91.17"
$ws.Range("E18").Value = "  -1.42%  "

$ws.Range("D19").Value = "'0.06652"
$ws.Range("E19").Value = "  +0.24%  "

$ws.Range("D20").Value = "'18.10"
$ws.Range("E20").Value = "  +1.19%  "

$ws.Range("E21").Value = "  +0.40%  "

$ws.Range("D22").Value = "'6.107"
$ws.Range("E22").Value = "  -1.99%  "

$ws.Range("D23").Value = "28.194.09"
$ws.Range("E23").Value = "  +0.18%  "

$ws.Range("D24").Value = "'11.21"
$ws.Range("E24").Value = "  -2.26%  "

$ws.Range("D25").Value = "'2.272"
$ws.Range("E25").Value = "  -1.78%  "

$ws.Range("E26").Value = "  +1.34%  "

$ws.Range("D27").Value = "2.098.40"
$ws.Range("E27").Value = "  -0.49%  "

$ws.Range("E28").Value = "  -0.85%  "

$ws.Range("D29").Value = "'156.34"
$ws.Range("E29").Value = "  -1.16%  "

$ws.Range("D30").Value = "'126.64"
$ws.Range("E30").Value = "  +0.04%  "

$ws.Range("E31").Value = "  -0.81%  "

$ws.Range("D32").Value = "'1.049"
$ws.Range("E32").Value = "  -3.11%  "

$ws.Range("E33").Value = "  +0.34%  "

$ws.Range("E34").Value = "  -0.12%  "

$ws.Range("D35").Value = "'9.723"
$ws.Range("E35").Value = "  +1.44%  "

$ws.Range("D36").Value = "'0.02454"
$ws.Range("E36").Value = "  +2.30%  "

$ws.Range("D37").Value = "'0.06549"
$ws.Range("E37").Value = "  -0.78%  "

$ws.Range("D38").Value = "'1.232"
$ws.Range("E38").Value = "  +0.48%  "

$ws.Range("D39").Value = "'0.2176"
$ws.Range("E39").Value = "  -1.40%  "

$ws.Range("D40").Value = "'1.241"
$ws.Range("E40").Value = "  -4.62%  "

$ws.Range("D41").Value = "'0.6379"
$ws.Range("E41").Value = "  -1.51%  "

$ws.Range("E42").Value = "  -0.26%  "

$ws.Range("D43").Value = "'4.898"
$ws.Range("E43").Value = "  -1.61%  "

$ws.Range("D44").Value = "'0.6057"
$ws.Range("E44").Value = "  -0.90%  "

$ws.Range("D45").Value = "'13.18"
$ws.Range("E45").Value = "  -1.22%  "

$ws.Range("D46").Value = "'1.298"
$ws.Range("E46").Value = "  -1.01%  "

$ws.Range("D47").Value = "'3.683"
$ws.Range("E47").Value = "  -0.12%  "

$ws.Range("D48").Value = "'2.005"
$ws.Range("E48").Value = "  +0.03%  "

$ws.Range("D49").Value = "'1.220"
$ws.Range("E49").Value = "  +2.19%  "

$ws.Range("D50").Value = "'120.87"
$ws.Range("E50").Value = "  -1.11%  "

$ws.Range("D51").Value = "'80.73"
$ws.Range("E51").Value = "  +2.28%  "
